# "se actualizan los CP 4 y 10"
# CP004_modo_oscuro (row 5) and CP010_crear_historia (row 11) gain their
# "Dato00x" values in the DataPrueba sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CP004_modo_oscuro (row 5): add Dato003 value ---
$ws.Range("D5").Value = "dark-mode"

# --- CP010_crear_historia (row 11): fill Dato001-Dato003 like the other CPs ---
$ws.Range("B11").Value = "jisola.tsoft@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:jisola.tsoft@gmail.com")
$ws.Range("B11").Style = "Hyperlink"
$ws.Range("C11").Value = 12061990
$ws.Range("D11").Value = "historia de prueba2"

# --- match the author's final selection position ---
[void]$ws.Range("D12").Select()
